$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D31").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"
$ws.Range("B32").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"
$ws.Range("D51").Value = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"
$ws.Range("B62").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
